$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 790, shifting rows 790:831 down to 791:832.
$ws.Rows(790).Insert()

# Populate the newly-inserted row 790 with the new data point.
# Column A holds a date-like string that must stay literal text (not be
# auto-converted to a date serial), so assign it with a leading apostrophe
# (Excel's "force text" convention) and then reset the cell style back to
# Normal so no stray quote-prefix formatting is left behind.
$ws.Range("A790").Value = "'2026/02/14"
$ws.Range("A790").Style = "Normal"

$ws.Range("B790").Value = "土"
$ws.Range("C790").Value = 1
$ws.Range("D790").Value = 24
